$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 399
$ws.Range("I6").Value = 400
$ws.Range("J6").Value = 398
$ws.Range("K6").Value = 1200
$ws.Range("L6").Value = 1194
$ws.Range("M6").Value = -1088
$ws.Range("N6").Value = -1418
$ws.Range("H64").Value = 2900
$ws.Range("I64").Value = 2752.6316
$ws.Range("J64").Value = 3100
$ws.Range("K64").Value = 2752.6316
$ws.Range("L64").Value = 3100
$ws.Range("M64").Value = -2504.6316
$ws.Range("N64").Value = -3596
$ws.Range("H67").Value = 2900
$ws.Range("I67").Value = 2752.6316
$ws.Range("J67").Value = 3100
$ws.Range("K67").Value = 2752.6316
$ws.Range("L67").Value = 3100
$ws.Range("M67").Value = -1894.6316
$ws.Range("N67").Value = -4816
$ws.Range("H74").Value = 7798.952
$ws.Range("I74").Value = 10657.5
$ws.Range("J74").Value = 3987.5557
$ws.Range("K74").Value = 10657.5
$ws.Range("L74").Value = 3987.5557
$ws.Range("M74").Value = -9721.5
$ws.Range("N74").Value = -5859.5557
$ws.Range("H77").Value = 7798.952
$ws.Range("I77").Value = 10657.5
$ws.Range("J77").Value = 3987.5557
$ws.Range("K77").Value = 53287.5
$ws.Range("L77").Value = 19937.7785
$ws.Range("M77").Value = -48607.5
$ws.Range("N77").Value = -29297.7785
$ws.Range("H111").Value = 20836118
$ws.Range("I111").Value = 27779826
$ws.Range("J111").Value = 4997.5
$ws.Range("K111").Value = 83339478
$ws.Range("L111").Value = 14992.5
$ws.Range("M111").Value = -83336411
$ws.Range("N111").Value = -21126.5
$ws.Range("H113").Value = 21741804
$ws.Range("I113").Value = 38463676
$ws.Range("J113").Value = 3368.5
$ws.Range("K113").Value = 38463676
$ws.Range("L113").Value = 3368.5
$ws.Range("M113").Value = -38460422
$ws.Range("N113").Value = -9876.5
$ws.Range("H125").Value = 2308.4443
$ws.Range("J125").Value = 2534.5
$ws.Range("L125").Value = 22810.5
$ws.Range("N125").Value = -27730.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1480.68
$ws.Range("I61").Value = 1310.3043
$ws.Range("K61").Value = 1310.3043
$ws.Range("M61").Value = -1098.3043
$ws.Range("H63").Value = 2831.5386
$ws.Range("I63").Value = 1758.5714
$ws.Range("J63").Value = 4083.3333
$ws.Range("K63").Value = 1758.5714
$ws.Range("L63").Value = 4083.3333
$ws.Range("M63").Value = -1072.5714
$ws.Range("N63").Value = -5455.3333
$ws.Range("H66").Value = 2831.5386
$ws.Range("I66").Value = 1758.5714
$ws.Range("J66").Value = 4083.3333
$ws.Range("K66").Value = 8792.857
$ws.Range("L66").Value = 20416.6665
$ws.Range("M66").Value = -5360.857
$ws.Range("N66").Value = -27280.6665
$ws.Range("H74").Value = 894.86365
$ws.Range("I74").Value = 875.5714
$ws.Range("J74").Value = 1300
$ws.Range("K74").Value = 875.5714
$ws.Range("L74").Value = 1300
$ws.Range("M74").Value = -1.57140000000004
$ws.Range("N74").Value = -3048
$ws.Range("H77").Value = 894.86365
$ws.Range("I77").Value = 875.5714
$ws.Range("J77").Value = 1300
$ws.Range("K77").Value = 4377.857
$ws.Range("L77").Value = 6500
$ws.Range("M77").Value = -9.856999999999971
$ws.Range("N77").Value = -15236
$ws.Range("H117").Value = 25848
$ws.Range("J117").Value = 25848
$ws.Range("L117").Value = 25848
$ws.Range("N117").Value = -35026
$ws.Range("H118").Value = 37818.5
$ws.Range("J118").Value = 37818.5
$ws.Range("L118").Value = 37818.5
$ws.Range("N118").Value = -41132.5
$ws.Range("H121").Value = 14142.917
$ws.Range("J121").Value = 14142.917
$ws.Range("L121").Value = 14142.917
$ws.Range("N121").Value = -17636.917
$ws.Range("H122").Value = 1610.5186
$ws.Range("I122").Value = 1512.3914
$ws.Range("J122").Value = 2174.75
$ws.Range("K122").Value = 4537.174199999999
$ws.Range("L122").Value = 6524.25
$ws.Range("M122").Value = -2087.174199999999
$ws.Range("N122").Value = -11424.25
$ws.Range("H128").Value = 56400
$ws.Range("J128").Value = 56400
$ws.Range("L128").Value = 56400
$ws.Range("N128").Value = -66360
$ws.Range("H136").Value = 1480.68
$ws.Range("I136").Value = 1310.3043
$ws.Range("K136").Value = 3930.9129
$ws.Range("M136").Value = -1380.9129

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2134.375
$ws.Range("I122").Value = 1252.2
$ws.Range("J122").Value = 3604.6667
$ws.Range("K122").Value = 3756.6
$ws.Range("L122").Value = 10814.0001
$ws.Range("M122").Value = -1306.6
$ws.Range("N122").Value = -15714.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I7").Value = 125
$ws.Range("J7").Value = 248.66667
$ws.Range("K7").Value = 375
$ws.Range("L7").Value = 746.00001
$ws.Range("M7").Value = -263
$ws.Range("N7").Value = -970.00001
$ws.Range("H16").Value = 4988.4
$ws.Range("I16").Value = 298.66666
$ws.Range("J16").Value = 6998.2856
$ws.Range("K16").Value = 895.9999799999999
$ws.Range("L16").Value = 20994.8568
$ws.Range("M16").Value = -722.9999799999999
$ws.Range("N16").Value = -21340.8568
$ws.Range("H92").Value = 375
$ws.Range("I92").Value = 375
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1125
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = 123
$ws.Range("H121").Value = 954008.9
$ws.Range("I121").Value = 981
$ws.Range("J121").Value = 1112846.9
$ws.Range("K121").Value = 2943
$ws.Range("L121").Value = 3338540.7
$ws.Range("M121").Value = -1633
$ws.Range("N121").Value = -3341160.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28281.818
$ws.Range("I70").Value = 34320
$ws.Range("J70").Value = 4800
$ws.Range("K70").Value = 34320
$ws.Range("L70").Value = 4800
$ws.Range("M70").Value = -34050
$ws.Range("N70").Value = -5340
$ws.Range("H73").Value = 28281.818
$ws.Range("I73").Value = 34320
$ws.Range("J73").Value = 4800
$ws.Range("K73").Value = 34320
$ws.Range("L73").Value = 4800
$ws.Range("M73").Value = -33384
$ws.Range("N73").Value = -6672
$ws.Range("H122").Value = 2082.5454
$ws.Range("I122").Value = 2333.3333
$ws.Range("J122").Value = 1781.6
$ws.Range("K122").Value = 6999.999899999999
$ws.Range("L122").Value = 5344.799999999999
$ws.Range("M122").Value = -4549.999899999999
$ws.Range("N122").Value = -10244.8
$ws.Range("H132").Value = 3698.158
$ws.Range("I132").Value = 3196
$ws.Range("J132").Value = 4150.1
$ws.Range("K132").Value = 9588
$ws.Range("L132").Value = 12450.3
$ws.Range("M132").Value = -7058
$ws.Range("N132").Value = -17510.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6251601
$ws.Range("I7").Value = 6668107.5
$ws.Range("K7").Value = 6668107.5
$ws.Range("M7").Value = -6667995.5
$ws.Range("H40").Value = 1349.3214
$ws.Range("I40").Value = 1306.7037
$ws.Range("K40").Value = 1306.7037
$ws.Range("M40").Value = -1170.7037
$ws.Range("H61").Value = 2699.889
$ws.Range("I61").Value = 1824.75
$ws.Range("K61").Value = 1824.75
$ws.Range("M61").Value = -1622.75
$ws.Range("H68").Value = 1407.9791
$ws.Range("I68").Value = 1541.1936
$ws.Range("J68").Value = 1165.0588
$ws.Range("K68").Value = 1541.1936
$ws.Range("L68").Value = 1165.0588
$ws.Range("M68").Value = -792.1936000000001
$ws.Range("N68").Value = -2663.0588
$ws.Range("H71").Value = 1407.9791
$ws.Range("I71").Value = 1541.1936
$ws.Range("J71").Value = 1165.0588
$ws.Range("K71").Value = 7705.968000000001
$ws.Range("L71").Value = 5825.294
$ws.Range("M71").Value = -3961.968000000001
$ws.Range("N71").Value = -13313.294
$ws.Range("H113").Value = 2699.889
$ws.Range("I113").Value = 1824.75
$ws.Range("K113").Value = 1824.75
$ws.Range("M113").Value = 345.25
$ws.Range("H118").Value = 33167
$ws.Range("J118").Value = 33167
$ws.Range("L118").Value = 33167
$ws.Range("N118").Value = -36481
$ws.Range("H126").Value = 6251601
$ws.Range("I126").Value = 6668107.5
$ws.Range("K126").Value = 20004322.5
$ws.Range("M126").Value = -20001852.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7739
$ws.Range("I62").Value = 3665
$ws.Range("J62").Value = 13850
$ws.Range("K62").Value = 3665
$ws.Range("L62").Value = 13850
$ws.Range("M62").Value = -3041
$ws.Range("N62").Value = -15098
$ws.Range("H65").Value = 7739
$ws.Range("I65").Value = 3665
$ws.Range("J65").Value = 13850
$ws.Range("K65").Value = 18325
$ws.Range("L65").Value = 69250
$ws.Range("M65").Value = -15205
$ws.Range("N65").Value = -75490
$ws.Range("H105").Value = 41269.57
$ws.Range("J105").Value = 41269.57
$ws.Range("L105").Value = 41269.57
$ws.Range("N105").Value = -48257.57
$ws.Range("H122").Value = 13596.223
$ws.Range("I122").Value = 22251.2
$ws.Range("J122").Value = 2777.5
$ws.Range("K122").Value = 66753.60000000001
$ws.Range("L122").Value = 8332.5
$ws.Range("M122").Value = -64303.60000000001
$ws.Range("N122").Value = -13232.5
$ws.Range("H125").Value = 47708.125
$ws.Range("J125").Value = 47708.125
$ws.Range("L125").Value = 47708.125
$ws.Range("N125").Value = -57548.125
$ws.Range("H126").Value = 2762.6428
$ws.Range("I126").Value = 1834.625
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 5503.875
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -3033.875
$ws.Range("N126").Value = -16940
$ws.Range("H128").Value = 47272.223
$ws.Range("J128").Value = 47272.223
$ws.Range("L128").Value = 47272.223
$ws.Range("N128").Value = -57232.223
$ws.Range("H136").Value = 701.82355
$ws.Range("I136").Value = 699.65216
$ws.Range("J136").Value = 706.36365
$ws.Range("K136").Value = 2098.95648
$ws.Range("L136").Value = 2119.09095
$ws.Range("M136").Value = 451.0435200000002
$ws.Range("N136").Value = -7219.09095
